$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.026.30'
$ws.Range('E2').Value = '  +4.99%  '
$ws.Range('D3').Value = '2.241.62'
$ws.Range('E3').Value = '  +4.69%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('E5').Value = '  +4.22%  '
$ws.Range('E6').Value = '  +2.81%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '75.05'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +9.92%  '
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.603'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.96'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.52%  '
$ws.Range('E11').Value = '  +4.07%  '
$ws.Range('E12').Value = '  +4.95%  '
$ws.Range('E13').Value = '  +2.19%  '
$ws.Range('D14').Value = '2.578.68'
$ws.Range('E14').Value = '  +4.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.56'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').Value = '2.234.98'
$ws.Range('E16').Value = '  +3.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.789'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.49%  '
$ws.Range('D18').Value = '42.915.36'
$ws.Range('E18').Value = '  +5.06%  '
$ws.Range('E19').Value = '  +6.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.10'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.99'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.19%  '
$ws.Range('E22').Value = '  +4.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +19.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '229.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.80'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.25%  '
$ws.Range('E27').Value = '  +6.72%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.96'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +31.91%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.33%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.14'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '171.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0802'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.28'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.24%  '
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.108'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.43'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0330'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +20.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.94'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +15.05%  '
$ws.Range('E40').Value = '  +4.59%  '
$ws.Range('E41').Value = '  +12.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.07%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '59.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.63%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '104.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.67'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.482'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +34.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0988'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.94%  '
$ws.Range('E48').Value = '  +13.28%  '
$ws.Range('E49').Value = '  +4.50%  '
$ws.Range('E50').Value = '  +4.58%  '
$ws.Range('E51').Value = '  +3.37%  '
